$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 2-5 hold weekly records for "Fruto del paraíso" at
# "Macroferia Regional de Talca". The update shifts the data
# (Fecha, Volumen, Precio mínimo, Precio máximo, Precio promedio
# ponderado, Precio $/Kg) up one row in a cycle: row3->row2,
# row4->row3, row5->row4, and the original row2 values wrap to row5.

$ws.Range("D2").Value = 44277
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("P2").Value = 550

$ws.Range("D3").Value = 44284
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 500

$ws.Range("D4").Value = 44291
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("P4").Value = 550

$ws.Range("D5").Value = 44280
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 500
